# Generate Report for Handback
# Update the "dafdef12-7a35-4b2d-9c12-246e16d27400" row (row 3) across the
# Overview / zh-cn / de-de sheets to reflect that the file has now been
# handed back (instead of still being "Ready for handoff" / stuck on a
# stale handback version).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-20 16:56:54"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).AutoFit()

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-20 16:57:01"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).AutoFit()
